$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.635.50'
$ws.Range("E2").Value = '  -1.21%  '
$ws.Range("D3").Value = '3.068.21'
$ws.Range("E3").Value = '  -3.60%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '586.87'
$ws.Range("E5").Value = '  -1.18%  '
$ws.Range("D6").Value = '154.32'
$ws.Range("E6").Value = '  +3.74%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = '3.066.87'
$ws.Range("E9").Value = '  -3.45%  '
$ws.Range("E10").Value = '  -4.59%  '
$ws.Range("E11").Value = '  -3.90%  '
$ws.Range("D12").Value = '0.447'
$ws.Range("E12").Value = '  -3.31%  '
$ws.Range("D13").Value = '36.66'
$ws.Range("E13").Value = '  -3.18%  '
$ws.Range("E14").Value = '  -4.80%  '
$ws.Range("E15").Value = '  -2.26%  '
$ws.Range("D16").Value = '3.574.04'
$ws.Range("E16").Value = '  -3.65%  '
$ws.Range("D17").Value = '63.603.41'
$ws.Range("E17").Value = '  -0.93%  '
$ws.Range("D18").Value = '7.11'
$ws.Range("E18").Value = '  -3.58%  '
$ws.Range("D19").Value = '3.069.25'
$ws.Range("E19").Value = '  -3.56%  '
$ws.Range("D20").Value = '469.48'
$ws.Range("E20").Value = '  -0.77%  '
$ws.Range("E21").Value = '  -2.62%  '
$ws.Range("E22").Value = '  -5.18%  '
$ws.Range("E23").Value = '  -3.13%  '
$ws.Range("E24").Value = '  -0.88%  '
$ws.Range("D25").Value = '80.26'
$ws.Range("E25").Value = '  -2.00%  '
$ws.Range("D26").Value = '12.74'
$ws.Range("E26").Value = '  -4.38%  '
$ws.Range("D27").Value = '10.42'
$ws.Range("E27").Value = '  +4.18%  '
$ws.Range("E28").Value = '  -0.24%  '
$ws.Range("D29").Value = '7.36'
$ws.Range("E29").Value = '  +1.44%  '
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.18%  '
$ws.Range("E31").Value = '  -3.43%  '
$ws.Range("D32").Value = '2.13'
$ws.Range("E32").Value = '  -5.75%  '
$ws.Range("E33").Value = '  -9.14%  '
$ws.Range("D34").Value = '26.98'
$ws.Range("E34").Value = '  -5.22%  '
$ws.Range("D35").Value = '0.0₃0818'
$ws.Range("E35").Value = '  -4.94%  '
$ws.Range("E36").Value = '  -2.81%  '
$ws.Range("E37").Value = '  -4.88%  '
$ws.Range("D38").Value = '3.23'
$ws.Range("E38").Value = '  -4.01%  '
$ws.Range("E39").Value = '  -5.20%  '
$ws.Range("D40").Value = '50.48'
$ws.Range("E40").Value = '  -2.07%  '
$ws.Range("D41").Value = '9.11'
$ws.Range("E41").Value = '  -3.50%  '
$ws.Range("D42").Value = '436.33'
$ws.Range("E42").Value = '  -6.73%  '
$ws.Range("D43").Value = '0.287'
$ws.Range("E43").Value = '  -3.20%  '
$ws.Range("D44").Value = '40.32'
$ws.Range("E44").Value = '  +2.19%  '
$ws.Range("E46").Value = '  -5.39%  '
$ws.Range("D47").Value = '2.796.57'
$ws.Range("E47").Value = '  -4.88%  '
$ws.Range("D48").Value = '130.08'
$ws.Range("E48").Value = '  -2.10%  '
$ws.Range("D50").Value = '25.01'
$ws.Range("E50").Value = '  +1.41%  '
$ws.Range("E51").Value = '  -2.98%  '

Write-Host "Updated 77 cells"
